$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers
$ws.Range("A1").Value = "テストNO"
$ws.Range("B1").Value = "テスト項目"
$ws.Range("C1").Value = "テスト概要"
$ws.Range("D1").Value = "手順"
$ws.Range("E1").Value = "クラス"
$ws.Range("F1").Value = "メソッド"
$ws.Range("G1").Value = "検証項目"
$ws.Range("H1").Value = "実施日"
$ws.Range("I1").Value = "テスト結果"

# Row 2
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "テスト"
$ws.Range("C2").Value = "テストのテスト１"
$ws.Range("D2").Value = "1`n2`n3"
$ws.Range("E2").Value = "dense.TestCase1"
$ws.Range("F2").Value = "test1"
$ws.Range("G2").Value = "aaaa`naaaa`naaaa"
$ws.Range("H2").Value = ""

# Clear leftover original row-3 cells before repositioning content
$ws.Range("B3").ClearContents()

# Row 3
$ws.Range("C3").Value = "テストのテスト２`nテストのテスト3"
$ws.Range("E3").Value = "selenium.SeleniumTest1"
$ws.Range("F3").Value = "test0"
$ws.Range("G3").Value = "bbb`nccc"
$ws.Range("H3").Value = ""

# Wrap text styling
$ws.Range("D2").WrapText = $true
$ws.Range("G2").WrapText = $true
$ws.Range("H2").WrapText = $true
$ws.Range("C3").WrapText = $true
$ws.Range("G3").WrapText = $true
$ws.Range("H3").WrapText = $true

# Row heights
$ws.Rows.Item(2).RowHeight = 51
$ws.Rows.Item(3).RowHeight = 68

# Selection
[void]$ws.Range("I2").Select()
